# popMart: fix case cancel payment
#
# Row 3 in the "Config" sheet was an accidental duplicate of row 2 (same
# email/password/card/address, but a different product URL in column E).
# This clears row 3 back to blank (keeping only the inherited "Hyperlink"
# cell style on A3/E3, same as row 2's A/E columns), removes its two
# hyperlinks, and leaves row 2 (and its two hyperlinks) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear row 3's contents -------------------------------------------
# Deleting the row and re-inserting a blank one in its place makes Excel
# inherit the blank formatting from row 2 above automatically, which is
# exactly the s="1" (Hyperlink style, no value) left on A3/E3 and no
# formatting at all on B3/C3/D3/F3 in the target layout.
$ws.Rows("3:3").Delete()
$ws.Rows("3:3").Insert()

# --- 2. Remove row 3's hyperlinks -----------------------------------------
# This engine's Range.Hyperlinks.Delete() drops every hyperlink on the
# worksheet (not just the target range), so row 2's two hyperlinks need to
# be re-created afterwards.
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("E3").Hyperlinks.Delete()

# --- 3. Re-create row 2's hyperlinks (A2 mailto, E2 product URL) ---------
$origA2 = $ws.Range("A2").Value2
$origE2 = $ws.Range("E2").Value2

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:basikle2005@gmail.com")
$ws.Range("A2").Value2 = $origA2
$ws.Range("A2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.yodobashi.com/product/100000001006404769", "", "", "https://www.yodobashi.com/product/100000001006404769")
$ws.Range("E2").Value2 = $origE2
$ws.Range("E2").Style = "Hyperlink"

# --- 4. Update the selected cell in the sheet view -------------------------
$ws.Range("F10").Select()
